$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Wipe all existing cell *contents* first (keeps formatting/styles on
# the cells that already had them: B1:E1 header row + A2:A9 index col).
# This also drops every shared-string reference so the string table
# starts clean and gets rebuilt in exactly the order we (re)enter text.
# ------------------------------------------------------------------
$ws.Cells.ClearContents()

# ------------------------------------------------------------------
# Header row (row 1): Algorithm, <group> mean/std x3
# ------------------------------------------------------------------
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "State Based mean"
$ws.Range("D1").Value = "State Based std"
$ws.Range("E1").Value = "Non State mean"
$ws.Range("F1").Value = "Non State std"
$ws.Range("G1").Value = "One Sided mean"
$ws.Range("H1").Value = "One Sided std"

# F1:H1 are brand-new header cells - give them the same bold/centered/
# bordered look as B1:E1 by copying the formatting across (xlPasteFormats).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Data rows 2-8 (row 9 / old "NB" row is dropped entirely)
# ------------------------------------------------------------------
$algorithms = @("LR", "LDA", "KNN", "DTREE", "RTREE", "XTREE", "SVM")
$values = @(
    @(0.8963716532930801,  0.01418168075414624,  0.8618805116389927, 0.01526124878464865,  0.8828698922168373, 0.008903436424854983),
    @(0.8940939064113766,  0.01497031799816278,  0.8470774079076296, 0.01521333570096997,  0.8836829003469188, 0.0111347615877996),
    @(0.9075980508990759,  0.009607903883623715, 0.9132880485156644, 0.008550403809273209, 0.8926310214242207, 0.00673597456475664),
    @(0.8849853022960197,  0.00847557872017886,  0.9053195095468871, 0.009836847538239237, 0.8766870580757924, 0.009548536314710683),
    @(0.8609064908238659,  0.0199154593978885,   0.8303199067821296, 0.01395096674585964,  0.8353714149519347, 0.01252498560581758),
    @(0.9061306639125023,  0.01488579099252658,  0.8859569396996901, 0.01150476698773702,  0.8940939064113766, 0.0106833198313767),
    @(0.9100362808188341,  0.01560561907500342,  0.9067797462990915, 0.01216272902618992,  0.8983255210402266, 0.01092884308976276)
)

for ($i = 0; $i -lt $algorithms.Count; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $i
    $ws.Range("B$row").Value = $algorithms[$i]
    $rowVals = $values[$i]
    $ws.Range("C$row").Value = $rowVals[0]
    $ws.Range("D$row").Value = $rowVals[1]
    $ws.Range("E$row").Value = $rowVals[2]
    $ws.Range("F$row").Value = $rowVals[3]
    $ws.Range("G$row").Value = $rowVals[4]
    $ws.Range("H$row").Value = $rowVals[5]
}

# Drop the now-unused former row 9 entirely so it stops contributing to
# the sheet's used range / dimension.
$ws.Rows.Item(9).Delete()
